$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 -- this shifts the existing rows 13..31
# (and their formatting) down to 14..32, growing the used range to A1:R32,
# matching the updated <dimension ref="A1:R32"/>.
$ws.Rows.Item(13).Insert()

# The row that used to be row 12 keeps its place as row 12, but gets a new
# weekly record: updated date (Fecha) and updated volume (Volumen).
$ws.Range("D12").Value = 44526
$ws.Range("J12").Value = 100

# The newly inserted (blank) row 13 is populated with the data that used to
# live in row 12 before the edit above (i.e. the old weekly record moves
# down one row, same as every other row below it).
$ws.Range("A13").Value = 7
$ws.Range("B13").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C13").Value = "Ñuble"
$ws.Range("D13").Value = 44519
$ws.Range("E13").Value = 16
$ws.Range("F13").Value = 100112026
$ws.Range("G13").Value = "Haba"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 80
$ws.Range("K13").Value = 6000
$ws.Range("L13").Value = 7000
$ws.Range("M13").Value = 6500
$ws.Range("N13").Value = "$/saco 25 kilos"
$ws.Range("O13").Value = "Provincia de Diguillín"
$ws.Range("P13").Value = 260
$ws.Range("Q13").Value = 25
$ws.Range("R13").Value = "Hortaliza"
